$wb = $excel.ActiveWorkbook

# Add the new worksheet after the last existing sheet (Sheet1), so it becomes the third/last tab.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "year2 dimension"

# Populate the small ~TFM_MIG table, column by column (top to bottom),
# matching the order the cells were originally authored in.
$ws.Range("B3").Value = "~TFM_MIG"
$ws.Range("B4").Value = "attribute"
$ws.Range("B5").Value = "CUM"

$ws.Range("C4").Value = "year2"
$ws.Range("C5").Value = "BOH-2030"

$ws.Range("D4").Value = "process"
$ws.Range("D5").Value = "*COA*"

$ws.Range("E4").Value = "value"
$ws.Range("E5").Value = "*.5"

# Make the new sheet the active one with the selection sitting just below the table,
# matching the saved selection state in the source file.
$ws.Range("E6").Select()
$ws.Activate()
